$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 163
$ws.Cells.Item(2, 8).Value = "kitchens"
$ws.Cells.Item(2, 9).Value = "distractor"
$ws.Cells.Item(2, 11).Value = "f"
$ws.Cells.Item(2, 12).Value = "stimuli/img_q577a.png"
$ws.Cells.Item(2, 13).Value = 81.26470588235294
$ws.Cells.Item(2, 14).Value = 59.08823529411764
$ws.Cells.Item(2, 15).Value = 70.17647058823529
$ws.Cells.Item(2, 16).Value = 34
$ws.Cells.Item(2, 17).Value = 8
$ws.Cells.Item(2, 18).Value = 8
$ws.Cells.Item(2, 19).Value = 8

# Row 3
$ws.Cells.Item(3, 6).Value = 164
$ws.Cells.Item(3, 8).Value = "bedrooms"
$ws.Cells.Item(3, 9).Value = "target"
$ws.Cells.Item(3, 11).Value = "j"
$ws.Cells.Item(3, 12).Value = "stimuli/img_qgbyn.png"
$ws.Cells.Item(3, 13).Value = 65.08108108108108
$ws.Cells.Item(3, 14).Value = 40.10810810810811
$ws.Cells.Item(3, 15).Value = 52.5945945945946
$ws.Cells.Item(3, 16).Value = 37
$ws.Cells.Item(3, 17).Value = 4
$ws.Cells.Item(3, 18).Value = 4
$ws.Cells.Item(3, 19).Value = 4

# Row 4
$ws.Cells.Item(4, 6).Value = 165
$ws.Cells.Item(4, 8).Value = "kitchens"
$ws.Cells.Item(4, 9).Value = "distractor"
$ws.Cells.Item(4, 11).Value = "f"
$ws.Cells.Item(4, 12).Value = "stimuli/img_uegbb.png"
$ws.Cells.Item(4, 13).Value = 78.80952380952381
$ws.Cells.Item(4, 14).Value = 61.52380952380953
$ws.Cells.Item(4, 15).Value = 70.16666666666667
$ws.Cells.Item(4, 16).Value = 42
$ws.Cells.Item(4, 17).Value = 8
$ws.Cells.Item(4, 18).Value = 8
$ws.Cells.Item(4, 19).Value = 8

# Row 5
$ws.Cells.Item(5, 6).Value = 166
$ws.Cells.Item(5, 8).Value = "bedrooms"
$ws.Cells.Item(5, 9).Value = "target"
$ws.Cells.Item(5, 11).Value = "j"
$ws.Cells.Item(5, 12).Value = "stimuli/img_th7xh.png"
$ws.Cells.Item(5, 13).Value = 82.35897435897436
$ws.Cells.Item(5, 14).Value = 65.53846153846153
$ws.Cells.Item(5, 15).Value = 73.94871794871796
$ws.Cells.Item(5, 16).Value = 39
$ws.Cells.Item(5, 17).Value = 8
$ws.Cells.Item(5, 18).Value = 8
$ws.Cells.Item(5, 19).Value = 8

# Row 6
$ws.Cells.Item(6, 6).Value = 167
$ws.Cells.Item(6, 8).Value = "bedrooms"
$ws.Cells.Item(6, 9).Value = "target"
$ws.Cells.Item(6, 11).Value = "j"
$ws.Cells.Item(6, 12).Value = "stimuli/img_oou46.png"
$ws.Cells.Item(6, 13).Value = 75.70270270270271
$ws.Cells.Item(6, 14).Value = 54.86486486486486
$ws.Cells.Item(6, 15).Value = 65.28378378378379
$ws.Cells.Item(6, 16).Value = 37
$ws.Cells.Item(6, 17).Value = 6
$ws.Cells.Item(6, 18).Value = 6
$ws.Cells.Item(6, 19).Value = 6

# Row 7
$ws.Cells.Item(7, 6).Value = 168
$ws.Cells.Item(7, 8).Value = "kitchens"
$ws.Cells.Item(7, 9).Value = "distractor"
$ws.Cells.Item(7, 11).Value = "f"
$ws.Cells.Item(7, 12).Value = "stimuli/img_cv6mf.png"
$ws.Cells.Item(7, 13).Value = 66.8
$ws.Cells.Item(7, 14).Value = 42.08
$ws.Cells.Item(7, 15).Value = 54.44
$ws.Cells.Item(7, 16).Value = 25
$ws.Cells.Item(7, 17).Value = 4
$ws.Cells.Item(7, 18).Value = 4
$ws.Cells.Item(7, 19).Value = 4

# Row 8
$ws.Cells.Item(8, 6).Value = 169
$ws.Cells.Item(8, 8).Value = "kitchens"
$ws.Cells.Item(8, 9).Value = "distractor"
$ws.Cells.Item(8, 11).Value = "f"
$ws.Cells.Item(8, 12).Value = "stimuli/img_2b8fp.png"
$ws.Cells.Item(8, 13).Value = 73.89189189189189
$ws.Cells.Item(8, 14).Value = 51.45945945945946
$ws.Cells.Item(8, 15).Value = 62.67567567567568
$ws.Cells.Item(8, 16).Value = 37
$ws.Cells.Item(8, 17).Value = 6
$ws.Cells.Item(8, 18).Value = 6
$ws.Cells.Item(8, 19).Value = 6

# Row 9
$ws.Cells.Item(9, 6).Value = 170
$ws.Cells.Item(9, 8).Value = "bedrooms"
$ws.Cells.Item(9, 9).Value = "target"
$ws.Cells.Item(9, 11).Value = "j"
$ws.Cells.Item(9, 12).Value = "stimuli/img_rvssl.png"
$ws.Cells.Item(9, 13).Value = 74.25
$ws.Cells.Item(9, 14).Value = 54.33333333333334
$ws.Cells.Item(9, 15).Value = 64.29166666666667
$ws.Cells.Item(9, 16).Value = 36
$ws.Cells.Item(9, 17).Value = 6
$ws.Cells.Item(9, 18).Value = 6
$ws.Cells.Item(9, 19).Value = 6

# Row 10
$ws.Cells.Item(10, 6).Value = 171
$ws.Cells.Item(10, 8).Value = "bedrooms"
$ws.Cells.Item(10, 9).Value = "target"
$ws.Cells.Item(10, 11).Value = "j"
$ws.Cells.Item(10, 12).Value = "stimuli/img_le8uf.png"
$ws.Cells.Item(10, 13).Value = 12.88888888888889
$ws.Cells.Item(10, 14).Value = 9.222222222222221
$ws.Cells.Item(10, 15).Value = 11.05555555555556
$ws.Cells.Item(10, 16).Value = 36
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = 1
$ws.Cells.Item(10, 19).Value = 1

# Row 11
$ws.Cells.Item(11, 6).Value = 172
$ws.Cells.Item(11, 8).Value = "bedrooms"
$ws.Cells.Item(11, 9).Value = "target"
$ws.Cells.Item(11, 11).Value = "j"
$ws.Cells.Item(11, 12).Value = "stimuli/img_v8dra.png"
$ws.Cells.Item(11, 13).Value = 61.77272727272727
$ws.Cells.Item(11, 14).Value = 38.79545454545455
$ws.Cells.Item(11, 15).Value = 50.28409090909091
$ws.Cells.Item(11, 16).Value = 44
$ws.Cells.Item(11, 17).Value = 3
$ws.Cells.Item(11, 18).Value = 3
$ws.Cells.Item(11, 19).Value = 3

# Row 12
$ws.Cells.Item(12, 6).Value = 173
$ws.Cells.Item(12, 8).Value = "bedrooms"
$ws.Cells.Item(12, 9).Value = "target"
$ws.Cells.Item(12, 11).Value = "j"
$ws.Cells.Item(12, 12).Value = "stimuli/img_bj2gr.png"
$ws.Cells.Item(12, 13).Value = 65.25
$ws.Cells.Item(12, 14).Value = 44.8
$ws.Cells.Item(12, 15).Value = 55.025
$ws.Cells.Item(12, 16).Value = 40
$ws.Cells.Item(12, 17).Value = 4
$ws.Cells.Item(12, 18).Value = 4
$ws.Cells.Item(12, 19).Value = 4

# Row 13
$ws.Cells.Item(13, 6).Value = 174
$ws.Cells.Item(13, 8).Value = "bedrooms"
$ws.Cells.Item(13, 9).Value = "target"
$ws.Cells.Item(13, 11).Value = "j"
$ws.Cells.Item(13, 12).Value = "stimuli/img_wyctg.png"
$ws.Cells.Item(13, 13).Value = 33.44736842105263
$ws.Cells.Item(13, 14).Value = 11.39473684210526
$ws.Cells.Item(13, 15).Value = 22.42105263157895
$ws.Cells.Item(13, 16).Value = 38
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = 1
$ws.Cells.Item(13, 19).Value = 1

# Row 14
$ws.Cells.Item(14, 6).Value = 175
$ws.Cells.Item(14, 8).Value = "bedrooms"
$ws.Cells.Item(14, 9).Value = "target"
$ws.Cells.Item(14, 11).Value = "j"
$ws.Cells.Item(14, 12).Value = "stimuli/img_uxxo0.png"
$ws.Cells.Item(14, 13).Value = 71.74418604651163
$ws.Cells.Item(14, 14).Value = 48.44186046511628
$ws.Cells.Item(14, 15).Value = 60.09302325581395
$ws.Cells.Item(14, 16).Value = 43
$ws.Cells.Item(14, 17).Value = 5
$ws.Cells.Item(14, 18).Value = 5
$ws.Cells.Item(14, 19).Value = 5

# Row 15
$ws.Cells.Item(15, 6).Value = 176
$ws.Cells.Item(15, 8).Value = "bedrooms"
$ws.Cells.Item(15, 9).Value = "target"
$ws.Cells.Item(15, 11).Value = "j"
$ws.Cells.Item(15, 12).Value = "stimuli/img_t2ioc.png"
$ws.Cells.Item(15, 13).Value = 88.1891891891892
$ws.Cells.Item(15, 14).Value = 74.05405405405405
$ws.Cells.Item(15, 15).Value = 81.12162162162161
$ws.Cells.Item(15, 16).Value = 37
$ws.Cells.Item(15, 17).Value = 10
$ws.Cells.Item(15, 18).Value = 10
$ws.Cells.Item(15, 19).Value = 10

# Row 16
$ws.Cells.Item(16, 6).Value = 177
$ws.Cells.Item(16, 8).Value = "bedrooms"
$ws.Cells.Item(16, 9).Value = "target"
$ws.Cells.Item(16, 11).Value = "j"
$ws.Cells.Item(16, 12).Value = "stimuli/img_a9acb.png"
$ws.Cells.Item(16, 13).Value = 77.11428571428571
$ws.Cells.Item(16, 14).Value = 58.42857142857143
$ws.Cells.Item(16, 15).Value = 67.77142857142857
$ws.Cells.Item(16, 16).Value = 35
$ws.Cells.Item(16, 17).Value = 7
$ws.Cells.Item(16, 18).Value = 7
$ws.Cells.Item(16, 19).Value = 7

# Row 17
$ws.Cells.Item(17, 6).Value = 178
$ws.Cells.Item(17, 8).Value = "living_rooms"
$ws.Cells.Item(17, 9).Value = "distractor"
$ws.Cells.Item(17, 11).Value = "f"
$ws.Cells.Item(17, 12).Value = "stimuli/img_pbsj1.png"
$ws.Cells.Item(17, 13).Value = 73.88636363636364
$ws.Cells.Item(17, 14).Value = 51.52272727272727
$ws.Cells.Item(17, 15).Value = 62.70454545454545
$ws.Cells.Item(17, 16).Value = 44
$ws.Cells.Item(17, 17).Value = 6
$ws.Cells.Item(17, 18).Value = 6
$ws.Cells.Item(17, 19).Value = 6

# Row 18
$ws.Cells.Item(18, 6).Value = 179
$ws.Cells.Item(18, 8).Value = "bedrooms"
$ws.Cells.Item(18, 9).Value = "target"
$ws.Cells.Item(18, 11).Value = "j"
$ws.Cells.Item(18, 12).Value = "stimuli/img_x0u5z.png"
$ws.Cells.Item(18, 13).Value = 92
$ws.Cells.Item(18, 14).Value = 78.16216216216216
$ws.Cells.Item(18, 15).Value = 85.08108108108108
$ws.Cells.Item(18, 16).Value = 37
$ws.Cells.Item(18, 17).Value = 10
$ws.Cells.Item(18, 18).Value = 10
$ws.Cells.Item(18, 19).Value = 10

# Row 19
$ws.Cells.Item(19, 6).Value = 180
$ws.Cells.Item(19, 8).Value = "bedrooms"
$ws.Cells.Item(19, 9).Value = "target"
$ws.Cells.Item(19, 11).Value = "j"
$ws.Cells.Item(19, 12).Value = "stimuli/img_h0hbk.png"
$ws.Cells.Item(19, 13).Value = 86.80952380952381
$ws.Cells.Item(19, 14).Value = 69.19047619047619
$ws.Cells.Item(19, 15).Value = 78
$ws.Cells.Item(19, 16).Value = 42
$ws.Cells.Item(19, 17).Value = 9
$ws.Cells.Item(19, 18).Value = 9
$ws.Cells.Item(19, 19).Value = 9

# Row 20
$ws.Cells.Item(20, 6).Value = 181
$ws.Cells.Item(20, 8).Value = "kitchens"
$ws.Cells.Item(20, 9).Value = "distractor"
$ws.Cells.Item(20, 11).Value = "f"
$ws.Cells.Item(20, 12).Value = "stimuli/img_a220l.png"
$ws.Cells.Item(20, 13).Value = 79.45945945945945
$ws.Cells.Item(20, 14).Value = 60.97297297297298
$ws.Cells.Item(20, 15).Value = 70.21621621621621
$ws.Cells.Item(20, 16).Value = 37
$ws.Cells.Item(20, 17).Value = 8
$ws.Cells.Item(20, 18).Value = 8
$ws.Cells.Item(20, 19).Value = 8

# Row 21
$ws.Cells.Item(21, 6).Value = 182
$ws.Cells.Item(21, 8).Value = "bedrooms"
$ws.Cells.Item(21, 9).Value = "target"
$ws.Cells.Item(21, 11).Value = "j"
$ws.Cells.Item(21, 12).Value = "stimuli/img_2js6m.png"
$ws.Cells.Item(21, 13).Value = 40.02777777777778
$ws.Cells.Item(21, 14).Value = 20.88888888888889
$ws.Cells.Item(21, 15).Value = 30.45833333333334
$ws.Cells.Item(21, 16).Value = 36
$ws.Cells.Item(21, 17).Value = 2
$ws.Cells.Item(21, 18).Value = 2
$ws.Cells.Item(21, 19).Value = 2

# Row 22
$ws.Cells.Item(22, 6).Value = 183
$ws.Cells.Item(22, 8).Value = "bedrooms"
$ws.Cells.Item(22, 9).Value = "target"
$ws.Cells.Item(22, 11).Value = "j"
$ws.Cells.Item(22, 12).Value = "stimuli/img_2pk6v.png"
$ws.Cells.Item(22, 13).Value = 85.08108108108108
$ws.Cells.Item(22, 14).Value = 66.16216216216216
$ws.Cells.Item(22, 15).Value = 75.62162162162161
$ws.Cells.Item(22, 16).Value = 37
$ws.Cells.Item(22, 17).Value = 9
$ws.Cells.Item(22, 18).Value = 9
$ws.Cells.Item(22, 19).Value = 9

# Row 23
$ws.Cells.Item(23, 6).Value = 184
$ws.Cells.Item(23, 8).Value = "bedrooms"
$ws.Cells.Item(23, 9).Value = "target"
$ws.Cells.Item(23, 11).Value = "j"
$ws.Cells.Item(23, 12).Value = "stimuli/img_71mhq.png"
$ws.Cells.Item(23, 13).Value = 69.34210526315789
$ws.Cells.Item(23, 14).Value = 47.02631578947368
$ws.Cells.Item(23, 15).Value = 58.18421052631579
$ws.Cells.Item(23, 16).Value = 38
$ws.Cells.Item(23, 17).Value = 5
$ws.Cells.Item(23, 18).Value = 5
$ws.Cells.Item(23, 19).Value = 5

# Row 24
$ws.Cells.Item(24, 6).Value = 185
$ws.Cells.Item(24, 8).Value = "bedrooms"
$ws.Cells.Item(24, 9).Value = "target"
$ws.Cells.Item(24, 11).Value = "j"
$ws.Cells.Item(24, 12).Value = "stimuli/img_okvvw.png"
$ws.Cells.Item(24, 13).Value = 50.58333333333334
$ws.Cells.Item(24, 14).Value = 32.11111111111111
$ws.Cells.Item(24, 15).Value = 41.34722222222223
$ws.Cells.Item(24, 16).Value = 36
$ws.Cells.Item(24, 17).Value = 2
$ws.Cells.Item(24, 18).Value = 2
$ws.Cells.Item(24, 19).Value = 2

# Row 25
$ws.Cells.Item(25, 6).Value = 186
$ws.Cells.Item(25, 8).Value = "bedrooms"
$ws.Cells.Item(25, 9).Value = "target"
$ws.Cells.Item(25, 11).Value = "j"
$ws.Cells.Item(25, 12).Value = "stimuli/img_ybbmx.png"
$ws.Cells.Item(25, 13).Value = 55.24324324324324
$ws.Cells.Item(25, 14).Value = 36.75675675675676
$ws.Cells.Item(25, 15).Value = 46
$ws.Cells.Item(25, 16).Value = 37
$ws.Cells.Item(25, 17).Value = 3
$ws.Cells.Item(25, 18).Value = 3
$ws.Cells.Item(25, 19).Value = 3

# Row 26
$ws.Cells.Item(26, 6).Value = 187
$ws.Cells.Item(26, 8).Value = "bedrooms"
$ws.Cells.Item(26, 9).Value = "target"
$ws.Cells.Item(26, 11).Value = "j"
$ws.Cells.Item(26, 12).Value = "stimuli/img_5m6x4.png"
$ws.Cells.Item(26, 13).Value = 80.23076923076923
$ws.Cells.Item(26, 14).Value = 58.41025641025641
$ws.Cells.Item(26, 15).Value = 69.32051282051282
$ws.Cells.Item(26, 16).Value = 39
$ws.Cells.Item(26, 17).Value = 7
$ws.Cells.Item(26, 18).Value = 7
$ws.Cells.Item(26, 19).Value = 7

# Row 27
$ws.Cells.Item(27, 6).Value = 188
$ws.Cells.Item(27, 8).Value = "bedrooms"
$ws.Cells.Item(27, 9).Value = "target"
$ws.Cells.Item(27, 11).Value = "j"
$ws.Cells.Item(27, 12).Value = "stimuli/img_fqgem.png"
$ws.Cells.Item(27, 13).Value = 80.75
$ws.Cells.Item(27, 14).Value = 61.475
$ws.Cells.Item(27, 15).Value = 71.1125
$ws.Cells.Item(27, 16).Value = 40
$ws.Cells.Item(27, 17).Value = 8
$ws.Cells.Item(27, 18).Value = 8
$ws.Cells.Item(27, 19).Value = 8
